# Update "想去人数" (interested-count) values in 江西-漫展信息.xlsx
# Sheet "展览" (Exhibitions): rows 6, 8, 12
# Sheet "全部类型" (All types): rows 7, 9, 13

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F6").Value = 5401
$wsExpo.Range("F8").Value = 5389
$wsExpo.Range("F12").Value = 12

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 5401
$wsAll.Range("F9").Value = 5389
$wsAll.Range("F13").Value = 12
